# Update gh-pages output (generated at 456a3b4)
# Applies "want-to-go" count bumps across sheets, and refreshes the
# "全部类型" (all-types) roll-up sheet: the 萤火虫动漫游戏嘉年华 event
# entry is gone (its slot is reused by the now-later rows shifting up),
# and a fresh 广州·AP动漫游戏嘉年华 entry is appended at the end of
# the shifted block.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (exhibitions) - simple counter (column F) bumps
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 840
$ws1.Cells.Item(6, 6).Value = 1007
$ws1.Cells.Item(13, 6).Value = 1181
$ws1.Cells.Item(14, 6).Value = 26252
$ws1.Cells.Item(15, 6).Value = 2853
$ws1.Cells.Item(17, 6).Value = 199
$ws1.Cells.Item(20, 6).Value = 245
$ws1.Cells.Item(21, 6).Value = 484
$ws1.Cells.Item(23, 6).Value = 190
$ws1.Cells.Item(24, 6).Value = 302
$ws1.Cells.Item(26, 6).Value = 613
$ws1.Cells.Item(27, 6).Value = 143
$ws1.Cells.Item(30, 6).Value = 38

# ---------------------------------------------------------------
# Sheet "演出" (performances) - simple counter (column F) bumps
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(6, 6).Value = 324
$ws2.Cells.Item(10, 6).Value = 4167
$ws2.Cells.Item(18, 6).Value = 40
$ws2.Cells.Item(21, 6).Value = 4186

# ---------------------------------------------------------------
# Sheet "全部类型" (all types roll-up)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# Plain counter (column F) bumps that keep their row position
$ws4.Cells.Item(11, 6).Value = 324
$ws4.Cells.Item(14, 6).Value = 1007
$ws4.Cells.Item(20, 6).Value = 1181
$ws4.Cells.Item(28, 6).Value = 199
$ws4.Cells.Item(34, 6).Value = 245
$ws4.Cells.Item(35, 6).Value = 484
$ws4.Cells.Item(37, 6).Value = 302
$ws4.Cells.Item(39, 6).Value = 613
$ws4.Cells.Item(40, 6).Value = 40
$ws4.Cells.Item(41, 6).Value = 143
$ws4.Cells.Item(45, 6).Value = 38

# Row 21 (萤火虫动漫游戏嘉年华) drops out of the roll-up; rows 22-27
# shift up into 21-26, and a brand-new row is written at 27
# (广州·AP动漫游戏嘉年华). Re-write rows 21-27 in place, each with
# the full final (post-shift) values for every column.

# Row 21 <- old row 22 (冰兔2024线下live「过去和未来」)
$ws4.Cells.Item(21, 2).Value = "2024-07-20"
$ws4.Cells.Item(21, 3).Value = "广州·冰兔2024线下live「过去和未来」"
$ws4.Cells.Item(21, 4).Value = "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
$ws4.Cells.Item(21, 5).Value = "2024.07.20 20:00-07.20 22:00"
$ws4.Cells.Item(21, 6).Value = 64
$ws4.Cells.Item(21, 7).Value = 198
$ws4.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87546"
$ws4.Cells.Item(21, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"

# Row 22 <- old row 23 (跨越二次元ACG神级动漫世界巡回演唱会)
$ws4.Cells.Item(22, 2).Value = "2024-07-20"
$ws4.Cells.Item(22, 3).Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$ws4.Cells.Item(22, 4).Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$ws4.Cells.Item(22, 5).Value = "2024.07.20 19:30-07.20 21:10"
$ws4.Cells.Item(22, 6).Value = 254
$ws4.Cells.Item(22, 7).Value = 280
$ws4.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$ws4.Cells.Item(22, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

# Row 23 <- old row 24 (昨日重现——唯美英文经典歌曲演唱会)
$ws4.Cells.Item(23, 2).Value = "2024-07-21"
$ws4.Cells.Item(23, 3).Value = "广州·昨日重现——唯美英文经典歌曲演唱会"
$ws4.Cells.Item(23, 4).Value = "东风中路299号 广州中山纪念堂"
$ws4.Cells.Item(23, 5).Value = "2024.07.21 19:30-07.21 21:30"
$ws4.Cells.Item(23, 6).Value = 2
$ws4.Cells.Item(23, 7).Value = 100
$ws4.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86802"
$ws4.Cells.Item(23, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/DR8AvmXe1716802703006.jpeg"

# Row 24 <- old row 25 (燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024)
$ws4.Cells.Item(24, 2).Value = "2024-07-21"
$ws4.Cells.Item(24, 3).Value = "广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
$ws4.Cells.Item(24, 4).Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$ws4.Cells.Item(24, 5).Value = "2024.07.21 14:30-07.21 16:00"
$ws4.Cells.Item(24, 6).Value = 174
$ws4.Cells.Item(24, 7).Value = 280
$ws4.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87034"
$ws4.Cells.Item(24, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png"

# Row 25 <- old row 26 (【早鸟8折】"浪漫古典Ⅱ"百年经典传世名曲烛光音乐会)
$ws4.Cells.Item(25, 2).Value = "2024-07-26"
$ws4.Cells.Item(25, 3).Value = "广州·【早鸟8折】“浪漫古典Ⅱ”百年经典传世名曲烛光音乐会 "
$ws4.Cells.Item(25, 4).Value = "广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）"
$ws4.Cells.Item(25, 5).Value = "2024.07.26 20:00-07.26 21:30"
$ws4.Cells.Item(25, 6).Value = 1
$ws4.Cells.Item(25, 7).Value = 144
$ws4.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87726"
$ws4.Cells.Item(25, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/A8vhVlhn1717575084179.png"

# Row 26 <- old row 27 (萨克斯王子安德鲁·杨——2024经典&流行音乐巡回演出)
$ws4.Cells.Item(26, 2).Value = "2024-07-26"
$ws4.Cells.Item(26, 3).Value = "广州·萨克斯王子安德鲁·杨——2024经典&流行音乐巡回演出"
$ws4.Cells.Item(26, 4).Value = "龙凤街道革新路124号太古仓码头5号仓 广州太空间Live House"
$ws4.Cells.Item(26, 5).Value = "2024.07.26 20:00-07.26 21:30"
$ws4.Cells.Item(26, 6).Value = 3
$ws4.Cells.Item(26, 7).Value = 280
$ws4.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86635"
$ws4.Cells.Item(26, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/rciNih361716802006584.jpeg"

# Row 27 <- brand-new entry (广州·AP动漫游戏嘉年华)
$ws4.Cells.Item(27, 2).Value = "2024-07-27"
$ws4.Cells.Item(27, 3).Value = "广州·AP动漫游戏嘉年华"
$ws4.Cells.Item(27, 4).Value = "新港东路630-638号 南丰国际会展中心"
$ws4.Cells.Item(27, 5).Value = "2024.07.27 09:00-07.28 17:00"
$ws4.Cells.Item(27, 6).Value = 2853
$ws4.Cells.Item(27, 7).Value = 80
$ws4.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87213"
$ws4.Cells.Item(27, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/3Z8rGZPP1718164976101.jpeg"
